# Refresh the cryptocurrency tracker sheet: overwrite each coin's
# Price (column D) and Volume(1h) (column E) cell with the latest
# scraped reading, as produced by this run's GitHub Actions scrape.
#
# Price cells are plain text (to preserve the feed's literal
# digit-grouping, e.g. "26.400.57"), so each is momentarily given a
# Text number format while its value is assigned (this stops Excel
# from "helpfully" re-parsing the string as a number and dropping
# meaningful trailing zeros, e.g. "61.70" -> 61.7) and then has that
# formatting cleared again so the cell's style is left untouched,
# matching the original file. Volume cells are already unambiguous
# text (they carry a "%" and padding spaces) so no such dance is
# needed for column E.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.400.57"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +0.45%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.617.64"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +1.39%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "212.81"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.08%  "

$ws.Range("E6").Value = "  -0.19%  "

$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("E8").Value = "  +0.09%  "

$ws.Range("E9").Value = "  +0.15%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.15"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +0.90%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0847"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -0.43%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.846.30"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +1.45%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.631.82"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +2.24%  "

$ws.Range("E14").Value = "  +0.06%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.508"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.01%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.81"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -0.24%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "236.22"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +9.28%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "26.408.31"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.50%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.81"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +4.42%  "

$ws.Range("E20").Value = "  +0.33%  "

$ws.Range("E21").Value = "  -0.03%  "

$ws.Range("E22").Value = "  -0.46%  "

$ws.Range("E23").Value = "  +1.13%  "

$ws.Range("E24").Value = "  +2.92%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.95"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +1.54%  "

$ws.Range("E27").Value = "  +0.97%  "

$ws.Range("E28").Value = "  +0.31%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.53"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +2.41%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0496"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +0.30%  "

$ws.Range("E31").Value = "  -0.18%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.519.56"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +6.03%  "

$ws.Range("E33").Value = "  +1.35%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.96"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -0.02%  "

$ws.Range("E35").Value = "  +3.42%  "

$ws.Range("E36").Value = "  +0.27%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.568"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +1.88%  "

$ws.Range("E38").Value = "  +0.39%  "

$ws.Range("E39").Value = "  +0.61%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.89"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +2.18%  "

$ws.Range("E41").Value = "  -0.01%  "

$ws.Range("E42").Value = "  +1.69%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.757.83"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +1.45%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.762"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +0.42%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "61.70"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +1.19%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.906"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.10%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "90.42"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +3.89%  "

$ws.Range("E48").Value = "  +1.83%  "

$ws.Range("E49").Value = "  +0.20%  "

$ws.Range("E50").Value = "  +0.93%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.50"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.71%  "
